# Bai 13 Cai dat VSM - "Dieu chinh ten bai va ten tep"
# Slide 1 (title slide): the subtitle shape's lesson-name line "Chuong 13. ..."
# becomes "Bai 13. ..." and the placeholder is made taller (and shifted up)
# to fit the now 3-line subtitle text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The subtitle placeholder ("Rectangle 3") is the 2nd shape on the slide.
$sh = $s.Shapes.Item(2)

# --- Resize / reposition the placeholder ------------------------------
# before: off  x=611560  y=3645024   ext cx=7920880 cy=792088
# after : off  x=611560  y=3429000   ext cx=7920880 cy=1512888
# (Left/Width stay the same; only Top/Height move.)
$sh.Top = 270.0
$sh.Height = 119.1251

# --- Split "Chuong 13. Cai dat mo hinh khong " into two runs -----------
# "Chuong 13" (9 chars) -> "Bai 13", keeping the remaining
# ". Cai dat mo hinh khong " text (and its run formatting) intact.
$tr = $sh.TextFrame.TextRange
$lead = $tr.Characters(1, 9)
$lead.Text = "Bài 13"
